# "pending run test cases with loops"
# The test-case loop on Hoja1 gains a new entry: row 4's TestCaseName
# ("NavigateAllTabsTestCase", a duplicate of row 2) is replaced with a new
# pending test case, "NavigateFooterLinksTestCase", and the active
# selection moves down to A6 ready for the next row to be filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A4").Value = "NavigateFooterLinksTestCase"

[void]$ws.Range("A6").Select()
